# honeywell_brent.xlsx save_data regen:
#  - column G ("K") is recomputed to use the strike-count (K) metric
#    instead of the old "Strike#" values
#  - row 8 additionally gets refreshed IP/I0/IF (H/I/J) values from the
#    regenerated std/mean + s_vals calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") per row, row -> value
$newK = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 0
    6  = 1
    7  = 3
    8  = 1
    9  = 1
    10 = 0
    11 = 0
    12 = 1
    13 = 0
    14 = 1
    15 = 2
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 2
    21 = 3
    22 = 1
    23 = 2
    24 = 1
    25 = 0
    26 = 0
    27 = 2
    28 = 1
    29 = 0
    30 = 1
    31 = 2
    32 = 1
    33 = 1
    34 = 1
    35 = 2
    36 = 1
    37 = 1
    38 = 3
    39 = 1
    40 = 3
    41 = 2
    42 = 1
    43 = 5
    44 = 7
    45 = 2
    46 = 5
    47 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}

# Row 8 also had IP (H), I0 (I) and IF (J) regenerated
$ws.Cells.Item(8, 8).Value = 3
$ws.Cells.Item(8, 9).Value = 6
$ws.Cells.Item(8, 10).Value = 8
